# Re-theme the Pandoc syntax-highlighting character styles from the
# dark 'zenburn' palette to the light 'pygments/tango' palette, and
# drop the dark-gray code-block shading (upgraded bootstrap theme).
$d = $word.ActiveDocument

function Set-TokStyle {
    param($doc, [string]$StyleName, $OleColor, [bool]$Bold, [bool]$Italic,
          [bool]$SetBold, [bool]$SetItalic)
    $st = $doc.Styles($StyleName)
    $st.Font.Color = $OleColor
    if ($SetBold)   { $st.Font.Bold = $Bold }
    if ($SetItalic) { $st.Font.Italic = $Italic }
    # Style-level shading (w:shd) mirrors the removed dark code-block
    # background; clear it via the object model too (best effort).
    $st.Font.Shading.Texture = 0
    $st.Font.Shading.BackgroundPatternColor = -16777216
    $st.Font.Shading.ForegroundPatternColor = -16777216
}

Set-TokStyle $d "KeywordTok" 2125824 $true $false $true $false
Set-TokStyle $d "DataTypeTok" 8336 $false $false $false $false
Set-TokStyle $d "DecValTok" 7381056 $false $false $false $false
Set-TokStyle $d "BaseNTok" 7381056 $false $false $false $false
Set-TokStyle $d "FloatTok" 7381056 $false $false $false $false
Set-TokStyle $d "ConstantTok" 136 $false $false $true $false
Set-TokStyle $d "CharTok" 10514496 $false $false $false $false
Set-TokStyle $d "SpecialCharTok" 10514496 $false $false $false $false
Set-TokStyle $d "StringTok" 10514496 $false $false $false $false
Set-TokStyle $d "VerbatimStringTok" 10514496 $false $false $false $false
Set-TokStyle $d "SpecialStringTok" 8939195 $false $false $false $false
Set-TokStyle $d "ImportTok" -16777216 $false $false $false $false
Set-TokStyle $d "CommentTok" 11575392 $false $true $false $true
Set-TokStyle $d "DocumentationTok" 2171322 $false $true $false $true
Set-TokStyle $d "AnnotationTok" 11575392 $true $true $true $true
Set-TokStyle $d "CommentVarTok" 11575392 $true $true $true $true
Set-TokStyle $d "OtherTok" 2125824 $false $false $false $false
Set-TokStyle $d "FunctionTok" 8267782 $false $false $false $false
Set-TokStyle $d "VariableTok" 8132377 $false $false $false $false
Set-TokStyle $d "ControlFlowTok" 2125824 $true $false $true $false
Set-TokStyle $d "OperatorTok" 6710886 $false $false $false $false
Set-TokStyle $d "BuiltInTok" -16777216 $false $false $false $false
Set-TokStyle $d "ExtensionTok" -16777216 $false $false $false $false
Set-TokStyle $d "PreprocessorTok" 31420 $false $false $true $false
Set-TokStyle $d "AttributeTok" 2723965 $false $false $false $false
Set-TokStyle $d "RegionMarkerTok" -16777216 $false $false $false $false
Set-TokStyle $d "InformationTok" 11575392 $true $true $true $true
Set-TokStyle $d "WarningTok" 11575392 $true $true $true $true
Set-TokStyle $d "AlertTok" 255 $true $false $true $false
Set-TokStyle $d "ErrorTok" 255 $true $false $true $false
Set-TokStyle $d "NormalTok" -16777216 $false $false $false $false

# Drop the dark code-block background from the SourceCode paragraph style.
$sc = $d.Styles("SourceCode")
$sc.ParagraphFormat.Shading.Texture = 0
$sc.ParagraphFormat.Shading.BackgroundPatternColor = -16777216
$sc.ParagraphFormat.Shading.ForegroundPatternColor = -16777216
